$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a flat "Conta / Nome / Saldo" export, sorted descending by
# Saldo, followed by a blank separator row and a "Filtros aplicados" note.
#
# This edit:
#   1) Removes 8 trailing rows with (now stale / out of scope) negative
#      balances, including the old entry for account 004216657 (JOAO)
#      which had balance -12316.76.
#   2) Removes the row for account 004346716 (TIAGO), which had the same
#      balance (1002.87) as account 004581652 (CINCO) right below it.
#   3) Re-adds account 004216657 (JOAO) with an updated balance of
#      12988.91, positioned to keep the Saldo column sorted descending
#      (between CLOVIS at 14952.59 and HEITOR at 11067.07).

# --- 1) delete the 8 trailing rows (bottom of the sheet, above the blank
#        separator row), from the bottom up so row numbers stay valid ---
$ws.Rows.Item(250).Delete()   # 004415557 FILIPE      -42165.72
$ws.Rows.Item(249).Delete()   # 001651617 MIRELLA      -14777.02
$ws.Rows.Item(248).Delete()   # 004216657 JOAO         -12316.76 (stale)
$ws.Rows.Item(247).Delete()   # 004512434 CAIO          -6344.44
$ws.Rows.Item(246).Delete()   # 004259659 BENTO         -5898.93
$ws.Rows.Item(245).Delete()   # 004254210 MARCO         -5358.44
$ws.Rows.Item(244).Delete()   # 004436055 MARCO         -2168.75
$ws.Rows.Item(243).Delete()   # 005064129 THIAGO          -706.70

# --- 2) delete the TIAGO row (still at its original row number: none of
#        the deletions above were above it) ---
$ws.Rows.Item(28).Delete()    # 004346716 TIAGO          1002.87

# --- 3) insert the refreshed JOAO row in sorted position ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "004216657"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = "JOAO"
$ws.Range("C16").Value = 12988.91

Write-Output ("Final used rows: " + $ws.UsedRange.Rows.Count)
